$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the first 5 data rows (old rows 2-6), shifting rows 7-21 up to rows 2-16
$ws.Range("2:6").EntireRow.Delete() | Out-Null

# Append the 15 newly recorded data rows (new rows 17-31)
$newRows = @(
    @(3.90800370544684, 1.824196169285198, -6.420193484572113),
    @(-0.09310042271849372, 3.190837242564208, -1.485911830526877),
    @(-5.086900020557672, -0.347581442587992, 4.233577517212423),
    @(-4.690330617414796, 3.666637846680965, 4.422437811158368),
    @(-3.853545032563759, 7.179606515853127, 0.02150726318360263),
    @(-0.654346335781087, 7.300580988816252, -5.563441323452309),
    @(3.714518678644319, -2.467024099305689, -5.610007460651495),
    @(2.151245980966297, 0.7463607152954532, -6.372515623686779),
    @(-1.10870781817722, -0.5732807598478937, 3.703908173764338),
    @(-3.764231722211575, 3.633465962331838, 3.590699133977193),
    @(-3.505413231302479, 9.71845419680489, 0.237572531231097),
    @(-2.849029684327336, 8.877369792083568, -7.541428113895723),
    @(3.086845004493447, -5.60899220659439, -5.2191632692931),
    @(3.54143344769713, -2.492278899325752, -5.184806862815476),
    @(-0.4025364195714203, -2.233978587095854, 0.2824505248356344)
)

$startRow = 17
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
